$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.910.06"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "1.832.30"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'245.86"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").Value = "'0.6903"
$ws.Range("E6").Value = "  -2.35%  "

$ws.Range("D7").Value = "'0.9993"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3055"
$ws.Range("E8").Value = "  -2.52%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07646"
$ws.Range("E9").Value = "  -3.11%  "

$ws.Range("D10").Value = "'23.51"
$ws.Range("E10").Value = "  -4.33%  "

$ws.Range("D11").Value = "'0.07815"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "1.829.99"
$ws.Range("E12").Value = "  -3.60%  "

$ws.Range("D13").Value = "'5.069"
$ws.Range("E13").Value = "  -2.59%  "

$ws.Range("D14").Value = "'90.40"
$ws.Range("E14").Value = "  -3.29%  "

$ws.Range("D15").Value = "'0.6781"
$ws.Range("E15").Value = "  -3.34%  "

$ws.Range("D16").Value = "'6.421"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("D17").Value = "'0.000008307"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").Value = "28.914.00"
$ws.Range("E18").Value = "  -2.61%  "

$ws.Range("D19").Value = "'242.62"
$ws.Range("E19").Value = "  -3.79%  "

$ws.Range("D20").Value = "2.081.44"
$ws.Range("E20").Value = "  -3.51%  "

$ws.Range("D21").Value = "'12.67"
$ws.Range("E21").Value = "  -3.40%  "

$ws.Range("D22").Value = "'0.9990"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").Value = "'7.442"
$ws.Range("E23").Value = "  -2.55%  "

$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'0.1466"
$ws.Range("E25").Value = "  -5.73%  "

$ws.Range("D26").Value = "'160.90"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("D27").Value = "'8.794"
$ws.Range("E27").Value = "  -2.30%  "

$ws.Range("E28").Value = "  -2.93%  "

$ws.Range("D29").Value = "'1.561"
$ws.Range("E29").Value = "  +3.95%  "

$ws.Range("D30").Value = "'4.217"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").Value = "'4.140"
$ws.Range("E31").Value = "  -2.71%  "

$ws.Range("E32").Value = "  -2.31%  "

$ws.Range("D33").Value = "'0.05119"
$ws.Range("E33").Value = "  -3.75%  "

$ws.Range("D34").Value = "'0.7544"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("D35").Value = "'1.840"

$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("D37").Value = "'2.675"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("E38").Value = "  -2.64%  "

$ws.Range("D39").Value = "1.227.40"
$ws.Range("E39").Value = "  -3.94%  "

$ws.Range("E40").Value = "  -2.90%  "

$ws.Range("D41").Value = "'0.9258"
$ws.Range("E41").Value = "  +3.44%  "

$ws.Range("D42").Value = "'108.91"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "'0.9986"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").Value = "'5.722"
$ws.Range("E44").Value = "  -6.21%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.537"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.980.54"
$ws.Range("E46").Value = "  -3.25%  "

$ws.Range("D47").Value = "'0.5167"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("E48").Value = "  -5.50%  "

$ws.Range("D49").Value = "'63.96"
$ws.Range("E49").Value = "  -10.64%  "

$ws.Range("E50").Value = "  -3.14%  "

$ws.Range("D51").Value = "'0.4189"
$ws.Range("E51").Value = "  -2.83%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
